# The "municipio-nombre" column (K) was previously curated/tagged as a
# measure; it is now re-tagged as a dimension, consistent with the
# neighbouring "provincia-nombre" (L) and "comarca-nombre" (M) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = "sdmx-dimension:refArea"
$ws.Range("K3").Value = "dim"
$ws.Range("K4").Value = "URI-Municipio"
